$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B4").Value = "inf"
$ws.Range("B6").Value = 913418.5637610762
$ws.Range("B7").Value = 2962731.063294502
$ws.Range("B8").Value = 19557809.19380879
$ws.Range("B10").Value = 6131791.906567112

# ---- Sheet: Costs and Revenues ----
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 594389.6782579747
$ws.Range("D2").Value = 594406.0775322419
$ws.Range("E2").Value = 189405.2487033564
$ws.Range("F2").Value = 194319.4617185486
$ws.Range("G2").Value = 190176.3122611334
$ws.Range("H2").Value = 189497.6383834903
$ws.Range("I2").Value = 183714.0373709966
$ws.Range("J2").Value = 186798.0279346411
$ws.Range("L2").Value = 188050.2058709903
$ws.Range("M2").Value = 193369.8844759364
$ws.Range("N2").Value = 186751.3439793044
$ws.Range("O2").Value = 187310.3736662395
$ws.Range("P2").Value = 183714.0373709966
$ws.Range("B3").Value = 288523.4171191893
$ws.Range("C3").Value = 40018.63863282887
$ws.Range("B4").Value = 438396.5295836807
$ws.Range("C4").Value = 423078.7314798642
$ws.Range("E4").Value = 12355.45755378897
$ws.Range("F4").Value = 17269.67056898116
$ws.Range("G4").Value = 13126.52111156592
$ws.Range("H4").Value = 12447.84723392284
$ws.Range("I4").Value = 6664.246221429168
$ws.Range("J4").Value = 9748.23678507367
$ws.Range("K4").Value = 9748.23678507367
$ws.Range("L4").Value = 11000.41472142284
$ws.Range("M4").Value = 16320.09332636897
$ws.Range("N4").Value = 9701.552829736969
$ws.Range("O4").Value = 10260.58251667206
$ws.Range("P4").Value = 6664.246221429168
$ws.Range("B5").Value = 40154.92008928722
$ws.Range("B6").Value = -172685.1885341826
$ws.Range("C6").Value = 90188.12307033759
$ws.Range("D6").Value = 109546.0339935039
$ws.Range("E6").Value = 89659.30818125037
$ws.Range("M6").Value = 166645.9502520794
$ws.Range("N6").Value = 166645.9502520794
$ws.Range("O6").Value = 166645.9502520794

# ---- Sheet: Installed Capacities ----
$ws = $wb.Worksheets.Item("Installed Capacities")
$ws.Range("B3").Value = 299.4183527195973

# ---- Sheet: Added Capacities ----
$ws = $wb.Worksheets.Item("Added Capacities")
$ws.Range("B3").Value = 299.4183527195973
$ws.Range("C3").Value = 43.94745303966922

# ---- Sheet: DG Dispatch ----
$ws = $wb.Worksheets.Item("DG Dispatch")
$ws.Range("H2").Value = 351.7691192106338
$ws.Range("I2").Value = 261.2248202355961
$ws.Range("J2").Value = 162.6177924993845
$ws.Range("K2").Value = 168.7672874155736
$ws.Range("L2").Value = 154.8567770823587
$ws.Range("M2").Value = 125.5236059864442
$ws.Range("N2").Value = 120.7388558752824
$ws.Range("O2").Value = 133.71009560408
$ws.Range("P2").Value = 165.0858607382365
$ws.Range("Q2").Value = 195.1368909901325
$ws.Range("R2").Value = 239.6824396098621
$ws.Range("S2").Value = 239.9606215534526
$ws.Range("T2").Value = 220.2452372048549
$ws.Range("U2").Value = 248.8053101871216
$ws.Range("G3").Value = 161.1251639608333
$ws.Range("H3").Value = 138.8146328796908
$ws.Range("I3").Value = 123.3775744401053
$ws.Range("J3").Value = 124.6531352647845
$ws.Range("K3").Value = 97.59607155587346
$ws.Range("L3").Value = 62.7989356036918
$ws.Range("M3").Value = 44.68802237877131
$ws.Range("N3").Value = 24.58528318678628
$ws.Range("O3").Value = 55.31542844358495
$ws.Range("P3").Value = 72.95571951604451
$ws.Range("Q3").Value = 122.5132338691918
$ws.Range("R3").Value = 173.0698009520273
$ws.Range("S3").Value = 210.0693686688361
$ws.Range("T3").Value = 230.5394686220007
$ws.Range("U3").Value = 249.6568043141749
$ws.Range("G4").Value = 169.3267675485344
$ws.Range("H4").Value = 168.3990210914259
$ws.Range("I4").Value = 170.5138353165135
$ws.Range("J4").Value = 147.557418132161
$ws.Range("K4").Value = 125.9512909308596
$ws.Range("L4").Value = 116.9952695469455
$ws.Range("M4").Value = 118.5417351554377
$ws.Range("N4").Value = 104.1153689225513
$ws.Range("O4").Value = 126.1753670254217
$ws.Range("P4").Value = 136.1245540434928
$ws.Range("Q4").Value = 176.3638936224054
$ws.Range("R4").Value = 225.2074709056173
$ws.Range("S4").Value = 245.3020430631099
$ws.Range("T4").Value = 218.4813791433803
$ws.Range("U4").Value = 291.2246192313701

# ---- Sheet: PV Dispatch ----
$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G2").Value = 1.203691870229536
$ws.Range("H2").Value = 12.32730936598824
$ws.Range("I2").Value = 46.40533082702423
$ws.Range("J2").Value = 102.1618428708942
$ws.Range("K2").Value = 153.1141197377105
$ws.Range("L2").Value = 189.9516048112476
$ws.Range("M2").Value = 211.3577601084422
$ws.Range("N2").Value = 214.7777496347318
$ws.Range("O2").Value = 202.8085386001369
$ws.Range("P2").Value = 173.0923955538452
$ws.Range("Q2").Value = 129.9851804512499
$ws.Range("R2").Value = 75.61140944330616
$ws.Range("S2").Value = 27.42912849285558
$ws.Range("T2").Value = 5.269161161929796
$ws.Range("U2").Value = 0.09629534961836286
$ws.Range("G3").Value = 0.6440319284912093
$ws.Range("H3").Value = 6.219992572533522
$ws.Range("I3").Value = 22.17390630989471
$ws.Range("J3").Value = 60.84689373521553
$ws.Range("K3").Value = 103.9970329441265
$ws.Range("L3").Value = 139.8368448243742
$ws.Range("M3").Value = 163.1830022321805
$ws.Range("N3").Value = 167.5019707350887
$ws.Range("O3").Value = 153.231579056415
$ws.Range("P3").Value = 122.9818513274135
$ws.Range("Q3").Value = 82.21011073161472
$ws.Range("R3").Value = 39.98647359527001
$ws.Range("S3").Value = 11.96261060158013
$ws.Range("T3").Value = 2.595900624050093
$ws.Range("U3").Value = 0.04237052161126378
$ws.Range("G4").Value = 0.5399347344123885
$ws.Range("H4").Value = 4.800510638684694
$ws.Range("I4").Value = 16.2373100130562
$ws.Range("J4").Value = 38.17338572295586
$ws.Range("K4").Value = 62.7305991435484
$ws.Range("L4").Value = 80.27356951436548
$ws.Range("M4").Value = 84.63722386793467
$ws.Range("N4").Value = 82.62473985785221
$ws.Range("O4").Value = 76.31732046039836
$ws.Range("P4").Value = 65.30265187838557
$ws.Range("Q4").Value = 45.21217126065919
$ws.Range("R4").Value = 24.27742905821521
$ws.Range("S4").Value = 9.409589871532258
$ws.Range("T4").Value = 2.306993865216568
$ws.Range("U4").Value = 0.02945098551340304
$ws.Range("O14").Value = 323.254833436289
$ws.Range("U14").Value = 0.1534843523671645
$ws.Range("R15").Value = 63.73410582691734
$ws.Range("U16").Value = 0.04694167948934294

# ---- Sheet: Fed-in Capacity ----
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("J11").Value = 43.96511994473918
$ws.Range("K11").Value = 77.83409231047497
$ws.Range("K14").Value = 77.83409231047497
$ws.Range("K17").Value = 77.83409231047497
$ws.Range("I18").Value = 2.632694258632081
$ws.Range("Q18").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("J20").Value = 43.96511994473918
$ws.Range("K20").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("Q32").Value = 55.20189757157522
$ws.Range("J33").Value = 32.78366918850629
$ws.Range("K33").Value = 0
$ws.Range("Q33").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("Q38").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("Q41").Value = 55.20189757157522
$ws.Range("K42").Value = 0
$ws.Range("R44").Value = 0
$ws.Range("L46").Value = 42.24342240080415
$ws.Range("M46").Value = 46.41251183645587
$ws.Range("N46").Value = 44.67260637956007
$ws.Range("P46").Value = 35.10170646165069

# ---- Sheet: Unmet Demand ----
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("G2").Value = 22.38652970403251
$ws.Range("G5").Value = 22.20985652598358
$ws.Range("J11").Value = 57.97960852900644
$ws.Range("K11").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("O14").Value = 13.26380076792788
$ws.Range("Q14").Value = 62.73788191118649
$ws.Range("R15").Value = 137.4865345694366
$ws.Range("K17").Value = 0
$ws.Range("I18").Value = 107.5759826808776
$ws.Range("Q18").Value = 73.68933668991343
$ws.Range("K19").Value = 88.69611276210766
$ws.Range("L19").Value = 69.3214679998587
$ws.Range("M19").Value = 68.27639578424916
$ws.Range("J20").Value = 57.97960852900644
$ws.Range("K20").Value = 77.83409231047497
$ws.Range("K21").Value = 35.83310349479973
$ws.Range("K24").Value = 35.83310349479973
$ws.Range("K26").Value = 77.83409231047497
$ws.Range("K32").Value = 77.83409231047497
$ws.Range("Q32").Value = 62.73788191118649
$ws.Range("J33").Value = 55.7330047947075
$ws.Range("K33").Value = 35.83310349479973
$ws.Range("Q33").Value = 73.68933668991343
$ws.Range("P35").Value = 62.28773471909744
$ws.Range("K38").Value = 77.83409231047497
$ws.Range("Q38").Value = 117.9397794827617
$ws.Range("K39").Value = 35.83310349479973
$ws.Range("Q41").Value = 62.73788191118649
$ws.Range("K42").Value = 35.83310349479973
$ws.Range("R44").Value = 194.7774559358693
$ws.Range("L46").Value = 27.07804559905455
$ws.Range("M46").Value = 21.86388394779328
$ws.Range("N46").Value = 10.37262059349192
$ws.Range("P46").Value = 62.24014880585811

# ---- Sheet: Household Surplus ----
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B2").Value = 1069901.420864355
$ws.Range("B5").Value = 353284.9052198305
$ws.Range("B6").Value = 367044.7016623687
$ws.Range("B7").Value = 355443.883181606
$ws.Range("B8").Value = 353543.5963242053
$ws.Range("B9").Value = 337349.5134892231
$ws.Range("B10").Value = 345984.6870674277
$ws.Range("B12").Value = 349490.7852892054
$ws.Range("B13").Value = 364385.8853830545
$ws.Range("B14").Value = 345853.9719924849
$ws.Range("B15").Value = 347419.2551159033
$ws.Range("B16").Value = 337349.5134892231
